$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
  2 = @(-0.03874681001824887, 0.5532001388111951, 0.7502571438541419, 0.866173853134659, 0.8856697285236464, 22)
  3 = @(0.7741633195380929, 0.9734393599335567, 2.377041849710642, 1.541765821942698, 1.366235212037819, 21)
  4 = @(0.6367776243981704, 1.301153388049756, 4.054070903386229, 2.013472349794312, 1.959748952894371, 20)
  5 = @(0.7664563126367501, 0.9574959425809361, 1.631796648714983, 1.277417961637843, 1.049933490516606, 19)
  6 = @(0.5932326787493827, 0.7851613914910194, 1.106570861342902, 1.051936719267325, 0.8938885739127869, 18)
  7 = @(0.3312469404061975, 0.547657782996088, 0.4036688942095937, 0.6353494268586333, 0.558852289159396, 17)
  8 = @(0.3086347685158412, 0.4826359502433059, 0.3365528937250792, 0.580131789962487, 0.5073302392072996, 16)
}

foreach ($row in $values.Keys) {
  $rowVals = $values[$row]
  $ws.Cells.Item($row, 2).Value = $rowVals[0]
  $ws.Cells.Item($row, 3).Value = $rowVals[1]
  $ws.Cells.Item($row, 4).Value = $rowVals[2]
  $ws.Cells.Item($row, 5).Value = $rowVals[3]
  $ws.Cells.Item($row, 6).Value = $rowVals[4]
  $ws.Cells.Item($row, 7).Value = $rowVals[5]
}
